$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$values = @{
    2  = 2
    3  = -1
    4  = 4
    5  = 4
    6  = 3
    7  = 1
    9  = 2
    10 = -3
    11 = 9
    12 = -6
    13 = 2
    14 = -1
    15 = 9
    16 = 5
    17 = -4
    18 = 5
    19 = -2
    20 = 4
    21 = 2
    22 = -1
    24 = 5
    25 = 1
    26 = -3
    27 = 4
    28 = -2
    29 = -2
    30 = -6
    32 = 2
    33 = 1
    34 = 4
    35 = 4
    36 = -1
    37 = 1
    39 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
